# FedExShipments.xlsx - "Changes of webdriver exception"
#
# Rows 2-26, column P hold FedEx tracking numbers used by the webdriver
# tests. Each one is replaced with a freshly generated tracking number.
#
# The tracking numbers are long purely-numeric strings, but the sheet
# stores them as TEXT (shared-string) cells with the default/general
# number format - not as numbers. Assigning a numeric-looking string via
# .Value/.Value2/.Formula directly would make Excel coerce it to a
# number (and forcing text via NumberFormat="@" would stamp a new cell
# style onto the cell, which the original cells never had). Writing the
# value as a self-referential text formula and then collapsing it back
# to a literal with Copy + PasteSpecial(values-only) yields a plain text
# cell with no style side effects, matching the original authoring.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$trackingNumbers = @(
    "320018594180",
    "320018594190",
    "320018594227",
    "320018594249",
    "320018594282",
    "320018594308",
    "320018594330",
    "320018594352",
    "320018594385",
    "320018594400",
    "320018594444",
    "320018594466",
    "320018594499",
    "320018594514",
    "320018594547",
    "320018594569",
    "320018594606",
    "320018594628",
    "320018594650",
    "320018594672",
    "320018594709",
    "320018594710",
    "320018594720",
    "320018594731",
    "320018594742"
)

$startRow = 2
for ($i = 0; $i -lt $trackingNumbers.Length; $i++) {
    $cell = $ws.Cells.Item($startRow + $i, 16)
    $cell.Formula = "=""" + $trackingNumbers[$i] + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$excel.CutCopyMode = $false
